$d = $word.ActiveDocument

function Get-ParagraphIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd("`r", "`n").Trim()
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# The "Bibliografia" section ends with the line "Estudos de caso: EPIA de
# origem.". Immediately after it, the page footer/boilerplate used to be
# reproduced as three extra paragraphs:
#   1) a blank paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#      pages. Original theme under Creative Commons Attribution"
# The site rebuild dropped that boilerplate, so remove those three
# paragraphs while leaving the anchor paragraph (and everything after the
# boilerplate, e.g. the trailing blank paragraph and the page break) intact.
$anchorIndex = Get-ParagraphIndexByText $d "Estudos de caso: EPIA de origem."

if ($anchorIndex -gt 0) {
    for ($n = 0; $n -lt 3; $n++) {
        $target = $d.Paragraphs.Item($anchorIndex + 1)
        $target.Range.Delete()
    }
}
